# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 09:52"

# 2. Insert Rusia's updated figures into row 21 (previously Israel's row),
#    and push the old Israel figures (which used to live in row 21) down
#    into row 22 (previously Rusia's row). Net effect: Rusia now appears
#    before Israel in the list, matching the new sort order, and Rusia's
#    numbers are refreshed.
$ws.Range("A21").Value = "Rusia"
$ws.Range("B21").Value = 10131
$ws.Range("C21").Value = 1459
$ws.Range("D21").Value = 698
$ws.Range("E21").Value = 9357
$ws.Range("F21").Value = 8
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = 76

$ws.Range("A22").Value = "Israel"
$ws.Range("B22").Value = 9404
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 801
$ws.Range("E22").Value = 8530
$ws.Range("F22").Value = 147
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 73

# 3. Update Rumania's figures (row 32)
$ws.Range("E32").Value = 4004
$ws.Range("G32").Value = 9
$ws.Range("H32").Value = 229

# 4. Update Sri Lanka's figures (row 113)
$ws.Range("D113").Value = 47
$ws.Range("E113").Value = 135

# 5. Update Islas Feroe's figures (row 115)
$ws.Range("F115").Value = 0
